$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$testName = "Create and Delete CitizenShip From Excel"
$browser  = "chrome"
$passed   = "PASSED"
$failed   = "FAILED"

# Status values (column B) for the newly appended rows 505-515
$statuses = @($failed, $failed, $failed, $failed, $failed, $failed, $failed, $passed, $failed, $passed, $passed)

$startRow = 505
for ($i = 0; $i -lt $statuses.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $testName
    $ws.Cells.Item($row, 2).Value = $statuses[$i]
    $ws.Cells.Item($row, 3).Value = $browser
}
